$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells: "_old" suffix -> "_FV2410", "_new" suffix -> "_FV2504" ---
$oldHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

$newHeaders = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $oldHeaders[$i]
}

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $newHeaders[$i]
}

# --- Turn the data range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# --- Freeze header row (pane split below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
